$excel.DisplayAlerts = $false

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets("BOM-DM0260(Production_1V2)")

# Remove the two unused blank sheets that shipped with the original export
$wb.Worksheets("Sheet2").Delete()
$wb.Worksheets("Sheet3").Delete()

# Update the J2 FFC/FPC connector row: the board now mates with the
# OV9282 camera module instead of the IMX378 one, so the internal part
# Name / Footprint / DesignItemId / LibRef fields need to change to the
# new component record. (Manufacturer part number FH26W-33S-0.3SHW(60),
# designator J2, description, and manufacturer Hirose stay the same.)
$ws.Range("A9").Value = "'AC_PY003-OV9282_CON"
$ws.Range("G9").Value = "'AC-PY003-OV9282"
$ws.Range("H9").Value = "'CMP-005-000036-1"
$ws.Range("I9").Value = "'CMP-005-000036-1"
